$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column I values to 3.01 for the listed rows
$rowsToSet = @(5, 13, 16, 22, 26, 28, 43, 51, 55, 61, 62, 67, 70)
foreach ($r in $rowsToSet) {
    $ws.Range("I$r").Value = 3.01
}

# Swap G35 and G36 values
$ws.Range("G35").Value = 0
$ws.Range("G36").Value = 1

# Update the sheet view: change selection (this also resets the scrolled
# top-left cell back to the sheet's natural top-left)
$ws.Range("M6").Select()
